$d = $word.ActiveDocument

$replacements = @(
    @("57÷8=", "52÷8="),
    @("19÷8=", "73÷6="),
    @("86÷7=", "79÷3="),
    @("28÷5=", "40÷7="),
    @("37÷5=", "48÷3="),
    @("14÷9=", "23÷7="),
    @("37÷9=", "81÷8="),
    @("36÷2=", "46÷2="),
    @("30÷9=", "87÷8="),
    @("59÷4=", "13÷2="),
    @("52÷2=", "72÷7="),
    @("41÷8=", "45÷6="),
    @("31÷2=", "34÷3="),
    @("24÷9=", "48÷3="),
    @("87÷9=", "16÷5="),
    @("76÷4=", "28÷9="),
    @("76÷9=", "58÷7="),
    @("62÷4=", "30÷4="),
    @("55÷9=", "52÷5="),
    @("20÷8=", "78÷2="),
    @("52÷6=", "58÷9="),
    @("40÷5=", "69÷4="),
    @("88÷6=", "40÷7="),
    @("91÷6=", "79÷2="),
    @("21÷4=", "62÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
